# Generate Report for Handback
# Updates the status of the "7bdc1e01..." row (which failed handback transform)
# and records the error detail explaining the handback/handoff file name mismatch,
# for both the zh-cn and de-de localization sheets (and the shared Overview status).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the shared "Status" text (shared string) for the 7bdc1e01 row.
# This is referenced by Overview!E3, Overview!F3, zh-cn!C3 and de-de!C3 -
# updating the shared string text on one of them updates all of them.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the Error Detail (column P) for the 7bdc1e01 row on each language sheet.
$wsZhCn.Range("P3").Value = "Handback file name: xi0p5knm.iyp is different with handoff file name: 7bdc1e01-83c3-46d4-a605-6ae6357893b1.d2fc46f5d105aa3959e0e8dff6574525f3023fed.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: xi0p5knm.iyp is different with handoff file name: 7bdc1e01-83c3-46d4-a605-6ae6357893b1.d2fc46f5d105aa3959e0e8dff6574525f3023fed.de-de."

# Widen the Error Detail column so the long message is readable.
# (ColumnWidth uses character units; Excel stores width = ColumnWidth + 5/6 in the
# underlying XML, so 39.1666... here yields the target stored width of 40.)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666666
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666666
